# Classifications.xlsx - add eBird and Catalogue of Life taxonomy hierarchy rows,
# and restyle the Flora of Bhutan name cell (A7) with the new "loading taxon names"
# font treatment.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Row 7 (Flora of Bhutan) - Name cell picks up a distinct font treatment
# (Arial Narrow, 11pt, not bold) separate from the rest of the Name column.
$ws1.Range("A7").Font.Name = "Arial Narrow"
$ws1.Range("A7").Font.Size = 11
$ws1.Range("A7").Font.Bold = $false

# The row now spans the full table width, like the header/body rows above it.
$ws1.Cells.Item(7, 3).Style = "Normal"
$ws1.Cells.Item(7, 4).Style = "Normal"
$ws1.Cells.Item(7, 5).Style = "Normal"

# New row 8: eBird Taxonomy Hierarchy (2010)
$ws1.Range("A8").Value = "eBird Taxonomy Hierarchy (2010)"

# New row 9: Catalogue of Life Taxonomy Hierarchy
$ws1.Range("A9").Value = "Catalogue of Life Taxonomy Hierarchy"
$ws1.Cells.Item(9, 2).Style = "Normal"

# Leave the selection where the edit ended up.
$ws1.Range("B13").Select()
